# Updated symbol list on Mon Dec 26 22:54:13 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates. Values are stored as text in this sheet, so force
# a text number format before assigning, otherwise Excel would auto-convert
# the numeric-looking strings to real numbers (losing trailing zeros, etc.).
$priceUpdates = @{
    "D2"  = "242.96"
    "D3"  = "23.06"
    "D4"  = "5.421"
    "D5"  = "0.05913"
    "D6"  = "3.438"
    "D7"  = "6.522"
    "D8"  = "0.8103"
    "D9"  = "0.9320"
    "D11" = "0.07416"
    "D12" = "0.03237"
    "D13" = "0.03087"
    "D14" = "0.09362"
    "D15" = "3.861"
    "D16" = "0.001575"
    "D18" = "0.01121"
    "D19" = "0.005939"
    "D20" = "0.001256"
    "D21" = "0.004904"
    "D22" = "0.00006810"
    "D23" = "3.567"
    "D24" = "2.141"
    "D25" = "0.3232"
    "D40" = "0.03956"
    "D41" = "0.006549"
    "D42" = "0.1074"
    "D43" = "0.003004"
    "D44" = "0.008776"
    "D45" = "0.00005238"
    "D47" = "0.6709"
    "D48" = "0.002395"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Column E (Volume(1h) label) updates - "Bestin24h"/"Worstin24h" markers moved.
$ws.Range("E18").Value = "17OneONEBestin24h"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
